# Updates the cryptos price/volume table with freshly scraped values.
# Note: several "Price" (column D) values look numeric (e.g. "0.618",
# "229.62") but must stay TEXT, matching how the sheet already stores them
# (and to avoid float round-off like 0.618 -> 0.61799999999999999). Excel's
# own rule for forcing a numeric-looking entry to be treated as text is a
# leading apostrophe (quote-prefix), so those assignments below use
# "'<value>" - the apostrophe itself is not part of the stored text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.578.11'
$ws.Range('E2').Value = '  +0.02%  '

$ws.Range('D3').Value = '2.210.60'

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '''229.62'
$ws.Range('E5').Value = '  -1.15%  '

$ws.Range('D6').Value = '''0.618'
$ws.Range('E6').Value = '  -3.17%  '

$ws.Range('D7').Value = '''59.69'
$ws.Range('E7').Value = '  -7.03%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '''0.401'
$ws.Range('E9').Value = '  -2.27%  '

$ws.Range('D10').Value = '''57.60'
$ws.Range('E10').Value = '  -3.06%  '

$ws.Range('D11').Value = '''0.0891'
$ws.Range('E11').Value = '  -1.52%  '

$ws.Range('E12').Value = '  -1.36%  '

$ws.Range('D13').Value = '2.539.23'
$ws.Range('E13').Value = '  -1.96%  '

$ws.Range('D14').Value = '''15.40'
$ws.Range('E14').Value = '  -5.31%  '

$ws.Range('D15').Value = '''22.28'

$ws.Range('D16').Value = '''5.65'
$ws.Range('E16').Value = '  -0.80%  '

$ws.Range('D17').Value = '''0.793'
$ws.Range('E17').Value = '  -4.90%  '

$ws.Range('D18').Value = '2.213.50'
$ws.Range('E18').Value = '  -2.17%  '

$ws.Range('D19').Value = '41.517.92'
$ws.Range('E19').Value = '  +0.21%  '

$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  -1.97%  '

$ws.Range('D21').Value = '''71.94'
$ws.Range('E21').Value = '  -2.68%  '

$ws.Range('D22').Value = '''6.06'
$ws.Range('E22').Value = '  -2.33%  '

$ws.Range('D23').Value = '''242.09'
$ws.Range('E23').Value = '  -3.81%  '

$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  -0.15%  '

$ws.Range('D25').Value = '''2.35'
$ws.Range('E25').Value = '  -1.82%  '

$ws.Range('E26').Value = '  -2.14%  '

$ws.Range('D27').Value = '''9.67'
$ws.Range('E27').Value = '  -1.96%  '

$ws.Range('D28').Value = '''168.55'
$ws.Range('E28').Value = '  -2.89%  '

$ws.Range('E29').Value = '  -4.80%  '

$ws.Range('D30').Value = '''19.74'
$ws.Range('E30').Value = '  -3.80%  '

$ws.Range('D31').Value = '''1.41'
$ws.Range('E31').Value = '  -3.59%  '

$ws.Range('D32').Value = '''2.54'
$ws.Range('E32').Value = '  -9.44%  '

$ws.Range('E33').Value = '  -3.11%  '

$ws.Range('D34').Value = '''4.95'
$ws.Range('E34').Value = '  -2.02%  '

$ws.Range('D35').Value = '''4.63'
$ws.Range('E35').Value = '  -2.60%  '

$ws.Range('E36').Value = '  +1.77%  '

$ws.Range('D37').Value = '''6.47'
$ws.Range('E37').Value = '  -8.25%  '

$ws.Range('D38').Value = '''2.36'
$ws.Range('E38').Value = '  -4.43%  '

$ws.Range('D39').Value = '''3.57'
$ws.Range('E39').Value = '  -6.96%  '

$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  +0.13%  '

$ws.Range('D41').Value = '''0.000232'
$ws.Range('E41').Value = '  -13.84%  '

$ws.Range('E42').Value = '  -1.90%  '

$ws.Range('E43').Value = '  -3.18%  '

$ws.Range('D44').Value = '''0.0956'
$ws.Range('E44').Value = '  +1.07%  '

$ws.Range('E45').Value = '  -2.73%  '

$ws.Range('D46').Value = '''97.09'
$ws.Range('E46').Value = '  -5.61%  '

$ws.Range('D47').Value = '1.463.54'
$ws.Range('E47').Value = '  -3.18%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''16.35'
$ws.Range('E48').Value = '  -8.68%  '

$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '''4.24'
$ws.Range('E49').Value = '  -13.82%  '

$ws.Range('E50').Value = '  -1.34%  '

$ws.Range('D51').Value = '''1.07'
$ws.Range('E51').Value = '  -5.08%  '
